$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dim = $ws.UsedRange
$lastRow = $dim.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Recorded By")
    if ($cell.Value2 -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
}
